# Auto-generated edit script: apply Titan_Profits Leve price/profit updates
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()

$ws.Range("H76").Value = 3971182.5
$ws.Range("I76").Value = 4276369.5
$ws.Range("J76").Value = 3750
$ws.Range("K76").Value = 4276369.5
$ws.Range("L76").Value = 3750
$ws.Range("M76").Value = -4276054.5
$ws.Range("N76").Value = -4380

$ws.Range("H79").Value = 3971182.5
$ws.Range("I79").Value = 4276369.5
$ws.Range("J79").Value = 3750
$ws.Range("K79").Value = 4276369.5
$ws.Range("L79").Value = 3750
$ws.Range("M79").Value = -4275277.5
$ws.Range("N79").Value = -5934

$ws.Range("H80").Value = 1208.3704
$ws.Range("I80").Value = 440.33334
$ws.Range("J80").Value = 1427.8096
$ws.Range("K80").Value = 1321.00002
$ws.Range("L80").Value = 4283.4288
$ws.Range("M80").Value = -323.0000199999999
$ws.Range("N80").Value = -6279.4288

$ws.Range("H83").Value = 1208.3704
$ws.Range("I83").Value = 440.33334
$ws.Range("J83").Value = 1427.8096
$ws.Range("K83").Value = 3963.00006
$ws.Range("L83").Value = 12850.2864
$ws.Range("M83").Value = 1028.99994
$ws.Range("N83").Value = -22834.2864

$ws.Range("H132").Value = 244892.48
$ws.Range("I132").Value = 264813.1
$ws.Range("K132").Value = 794439.2999999999
$ws.Range("M132").Value = -791909.2999999999

$ws.Range("H138").Value = 7079113
$ws.Range("I138").Value = 3477457.2
$ws.Range("J138").Value = 8133256.5
$ws.Range("K138").Value = 10432371.6
$ws.Range("L138").Value = 24399769.5
$ws.Range("M138").Value = -10427231.6
$ws.Range("N138").Value = -24410049.5

$ws.Range("H139").Value = 49413.332
$ws.Range("J139").Value = 49413.332
$ws.Range("L139").Value = 49413.332
$ws.Range("N139").Value = -59693.332

$ws.Range("H141").Value = 1972.0646
$ws.Range("I141").Value = 1894.2759
$ws.Range("J141").Value = 3100
$ws.Range("K141").Value = 5682.8277
$ws.Range("L141").Value = 9300
$ws.Range("M141").Value = -502.8276999999998
$ws.Range("N141").Value = -19660

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H112").Value = 14721.5
$ws.Range("J112").Value = 14721.5
$ws.Range("L112").Value = 14721.5
$ws.Range("N112").Value = -17675.5

$ws.Range("H114").Value = 21699.334
$ws.Range("J114").Value = 21699.334
$ws.Range("L114").Value = 21699.334
$ws.Range("N114").Value = -30377.334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()

$ws.Range("H86").Value = 7728.5293
$ws.Range("I86").Value = 2919.6
$ws.Range("J86").Value = 9732.25
$ws.Range("K86").Value = 2919.6
$ws.Range("L86").Value = 9732.25
$ws.Range("M86").Value = -1796.6
$ws.Range("N86").Value = -11978.25

$ws.Range("H89").Value = 7728.5293
$ws.Range("I89").Value = 2919.6
$ws.Range("J89").Value = 9732.25
$ws.Range("K89").Value = 14598
$ws.Range("L89").Value = 48661.25
$ws.Range("M89").Value = -8982
$ws.Range("N89").Value = -59893.25

$ws.Range("H134").Value = 23258966
$ws.Range("I134").Value = 38463456
$ws.Range("K134").Value = 115390368
$ws.Range("M134").Value = -115387833

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1493.6818
$ws.Range("I31").Value = 961.1875
$ws.Range("J31").Value = 2913.6667
$ws.Range("K31").Value = 961.1875
$ws.Range("L31").Value = 2913.6667
$ws.Range("M31").Value = -666.1875
$ws.Range("N31").Value = -3503.6667

$ws.Range("H34").Value = 1493.6818
$ws.Range("I34").Value = 961.1875
$ws.Range("J34").Value = 2913.6667
$ws.Range("K34").Value = 961.1875
$ws.Range("L34").Value = 2913.6667
$ws.Range("M34").Value = -759.1875
$ws.Range("N34").Value = -3317.6667

$ws.Range("H122").Value = 1552.8
$ws.Range("I122").Value = 953.2778
$ws.Range("J122").Value = 2452.0833
$ws.Range("K122").Value = 2859.8334
$ws.Range("L122").Value = 7356.249899999999
$ws.Range("M122").Value = -409.8334
$ws.Range("N122").Value = -12256.2499

$ws.Range("H132").Value = 3285.5334
$ws.Range("I132").Value = 2330.375
$ws.Range("K132").Value = 6991.125
$ws.Range("M132").Value = -4461.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1219
$ws.Range("I5").Value = 478.41177
$ws.Range("J5").Value = 1959.5883
$ws.Range("K5").Value = 1435.23531
$ws.Range("L5").Value = 5878.7649
$ws.Range("M5").Value = -1323.23531
$ws.Range("N5").Value = -6102.7649

$ws.Range("H135").Value = 1219
$ws.Range("I135").Value = 478.41177
$ws.Range("J135").Value = 1959.5883
$ws.Range("K135").Value = 4305.70593
$ws.Range("L135").Value = 17636.2947
$ws.Range("M135").Value = -1770.70593
$ws.Range("N135").Value = -22706.2947

$ws.Range("H136").Value = 2510
$ws.Range("I136").Value = 1716.6666
$ws.Range("J136").Value = 2807.5
$ws.Range("K136").Value = 5149.9998
$ws.Range("L136").Value = 8422.5
$ws.Range("M136").Value = -49.9997999999996
$ws.Range("N136").Value = -18622.5

$ws.Range("H139").Value = 1848.1562
$ws.Range("I139").Value = 1801.3226
$ws.Range("J139").Value = 3300
$ws.Range("K139").Value = 5403.9678
$ws.Range("L139").Value = 9900
$ws.Range("M139").Value = -263.9678000000004
$ws.Range("N139").Value = -20180

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H54").Value = 50000
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()

$ws.Range("H80").Value = 2962.625
$ws.Range("I80").Value = 2850
$ws.Range("J80").Value = 3000.1667
$ws.Range("K80").Value = 2850
$ws.Range("L80").Value = 3000.1667
$ws.Range("M80").Value = -1852
$ws.Range("N80").Value = -4996.1667

$ws.Range("H83").Value = 2962.625
$ws.Range("I83").Value = 2850
$ws.Range("J83").Value = 3000.1667
$ws.Range("K83").Value = 14250
$ws.Range("L83").Value = 15000.8335
$ws.Range("M83").Value = -9258
$ws.Range("N83").Value = -24984.8335

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("M14").ClearContents()
$ws.Range("N14").ClearContents()

$ws.Range("H46").Value = 2345.4546
$ws.Range("I46").Value = 1280
$ws.Range("K46").Value = 1280
$ws.Range("M46").Value = -1092

$ws.Range("H55").Value = 530
$ws.Range("I55").Value = 477.55554
$ws.Range("J55").Value = 1002
$ws.Range("K55").Value = 477.55554
$ws.Range("L55").Value = 1002
$ws.Range("M55").Value = -304.55554
$ws.Range("N55").Value = -1348

$ws.Range("H122").Value = 3952.8235
$ws.Range("I122").Value = 3799.5
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 11398.5
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -8948.5
$ws.Range("N122").Value = -16900

$ws.Range("H130").Value = 10429
$ws.Range("J130").Value = 10429
$ws.Range("L130").Value = 10429
$ws.Range("N130").Value = -20469

$ws.Range("H133").Value = 50811.5
$ws.Range("J133").Value = 50811.5
$ws.Range("L133").Value = 50811.5
$ws.Range("N133").Value = -55871.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 2000
$ws.Range("J8").Value = 2000
$ws.Range("L8").Value = 2000
$ws.Range("N8").Value = -2280

$ws.Range("H11").Value = 13121
$ws.Range("I11").Value = 25052.5
$ws.Range("J11").Value = 5166.6665
$ws.Range("K11").Value = 25052.5
$ws.Range("L11").Value = 5166.6665
$ws.Range("M11").Value = -24910.5
$ws.Range("N11").Value = -5450.6665

$ws.Range("H123").Value = 29427.303
$ws.Range("J123").Value = 29427.303
$ws.Range("L123").Value = 29427.303
$ws.Range("N123").Value = -39227.303
